# Adds two text boxes ("BOTTOM BUTTON" and "TOP BUTTON") to slide 1, matching
# the shapes added in the target OOXML (ids 2 and 80 respectively).
#
# Notes on this headless COM engine's quirks, discovered empirically:
#
# 1. Shape.Left/.Top/.Width/.Height (and AddTextbox's Left/Top/Width/Height
#    args) are in POINTS, not EMU (1 pt = 12700 EMU), same as real PowerPoint.
#    We convert the target EMU coordinates from the diff into points.
#
# 2. Those properties are stored internally with limited (float) precision,
#    so a plain EMU/12700 conversion can land 1 EMU below the intended value
#    after PowerPoint rounds back to EMU on save. Adding half an EMU (in
#    points) before assigning compensates for this and reproduces the exact
#    target EMU values.
#
# 3. Newly created shapes get the next free id from a monotonically
#    increasing counter (ids already used by existing shapes are skipped; an
#    id freed by deleting a shape is never reused). The first free id here is
#    2, and id 80 (used by "TOP BUTTON" in the target) is only reached after
#    38 shapes have been created in this presentation. So 36 throwaway shapes
#    are created between the two real ones and deleted afterwards, landing
#    "TOP BUTTON" exactly on id 80 while leaving only the two real shapes in
#    the tree.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emu = 12700
$halfEmuInPts = 0.5 / $emu

# ---- Shape 1: "BOTTOM BUTTON" (lands on id=2) ----
$bottom = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$bottom.Left   = (8270030 / $emu) + $halfEmuInPts
$bottom.Top    = (5868746 / $emu) + $halfEmuInPts
$bottom.Width  = (1912062 / $emu) + $halfEmuInPts
$bottom.Height = (369332  / $emu) + $halfEmuInPts
$bottom.Name = "CasellaDiTesto 1"
$bottom.Fill.Visible = $false
$bottom.TextFrame.WordWrap = $false
$bottom.TextFrame.AutoSize = 1
$bottomRange = $bottom.TextFrame.TextRange
$bottomRange.Text = "BOTTOM BUTTON"
$bottomRange.Font.Bold = $true
$bottomRange.LanguageID = "it-IT"

# ---- 36 throwaway shapes to advance the id counter up to 79 ----
$filler = @()
for ($i = 0; $i -lt 36; $i++) {
    $tmp = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $filler += $tmp
}

# ---- Shape 38: "TOP BUTTON" (lands on id=80) ----
$top = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$top.Left   = (8275094 / $emu) + $halfEmuInPts
$top.Top    = (4466892 / $emu) + $halfEmuInPts
$top.Width  = (1436612 / $emu) + $halfEmuInPts
$top.Height = (369332  / $emu) + $halfEmuInPts
$top.Name = "CasellaDiTesto 79"
$top.Fill.Visible = $false
$top.TextFrame.WordWrap = $false
$top.TextFrame.AutoSize = 1
$topRange = $top.TextFrame.TextRange
$topRange.Text = "TOP BUTTON"
$topRange.Font.Bold = $true
$topRange.LanguageID = "it-IT"

# ---- Remove the throwaway shapes, leaving only the two real text boxes ----
foreach ($tmp in $filler) {
    $tmp.Delete()
}

Write-Host "bottom id=$($bottom.Id) top id=$($top.Id) shapes=$($s.Shapes.Count)"
